$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style source cells already present in the sheet:
#   A13 -> style "0" (Arial 10, no fill)
#   B1  -> style "1" (Cambria 11, green fill)
#   C1  -> style "2" (Cambria 11, no fill)
$styleA = $ws.Range("A13")
$style1 = $ws.Range("B1")
$style2 = $ws.Range("C1")

function Set-WithStyle($addr, $value, $src) {
    $rng = $ws.Range($addr)
    $rng.Value = $value
    $src.Copy()
    $rng.PasteSpecial(-4122)
}

# --- New "T" column (rows 2-11): per-row winning label, mirroring column A ---
Set-WithStyle "T2"  "Badewanne"      $style1
Set-WithStyle "T3"  "Dusche"         $style1
Set-WithStyle "T4"  "Handtuch"       $style1
Set-WithStyle "T5"  "Klobürste"      $style2
Set-WithStyle "T6"  "Klopapierrolle" $style1
Set-WithStyle "T7"  "Seife"          $style1
Set-WithStyle "T8"  "Shampoo"        $style1
Set-WithStyle "T9"  "Spiegel"        $style1
Set-WithStyle "T10" "Toilette"       $style1
Set-WithStyle "T11" "Waschbecken"    $style1

# --- A12: row label for the SUM row ---
Set-WithStyle "A12" "Summe" $styleA

# --- Row 17: "Rangfolge" (ranking) table ---
Set-WithStyle "A17" "Rangfolge" $styleA
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 9
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 3
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 7
$ws.Range("P17").Value = 8

# --- Rows 21-32: small validation/reference list ---
Set-WithStyle "A21" "Bad"          $styleA
Set-WithStyle "A22" "Toilette"     $styleA
Set-WithStyle "A23" "Klopapier"    $styleA
Set-WithStyle "A24" "Dusche"       $styleA
Set-WithStyle "A26" "Badewanne"    $styleA
Set-WithStyle "A27" "Waschbecken"  $styleA
Set-WithStyle "A28" "Handtuch"     $styleA
Set-WithStyle "A29" "Spiegel"      $styleA
Set-WithStyle "A30" "Shampoo"      $styleA
Set-WithStyle "A31" "Zahnbürste"   $styleA
Set-WithStyle "A32" "Zahnpasta"    $styleA

# --- Rows 39-48: object -> rank validation table ---
Set-WithStyle "A39" "Toilette"       $style1
$ws.Range("B39").Value = 1
Set-WithStyle "A40" "Klobürste"      $style2
$ws.Range("B40").Value = 2
Set-WithStyle "A41" "Seife"          $style1
$ws.Range("B41").Value = 3
Set-WithStyle "A42" "Waschbecken"    $style1
$ws.Range("B42").Value = 4
Set-WithStyle "A43" "Spiegel"        $style1
$ws.Range("B43").Value = 5
Set-WithStyle "A44" "Badewanne"      $style1
$ws.Range("B44").Value = 6
Set-WithStyle "A45" "Dusche"         $style1
$ws.Range("B45").Value = 7
Set-WithStyle "A46" "Klopapierrolle" $style1
$ws.Range("B46").Value = 8
Set-WithStyle "A47" "Handtuch"       $style1
$ws.Range("B47").Value = 9
Set-WithStyle "A48" "Shampoo"        $style1
$ws.Range("B48").Value = 10

# Restore selection to match the final workbook state.
[void]$ws.Range("F23").Select()
